$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.805.93'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.415.93'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.53%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.19'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.37'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.80%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.415.50'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.41%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.570'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -6.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.24'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.93%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.428'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.002.61'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +0.51%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.17'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -2.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000174'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -6.54%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '63.865.74'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.43%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.391.30'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.13'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -3.18%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.66'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.91%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '383.08'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +2.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.78'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.08%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.78%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.519'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000116'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.85%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.67'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.18%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.32%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.11'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +0.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.40'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.33%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.60%  '
$ws.Range('B33').Value = 'EthereumClassic'
$ws.Range('C33').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.01'
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.00'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.71%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -3.69%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '161.02'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.78%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +8.92%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.83'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.37%  '
$ws.Range('B40').Value = 'EnergySwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '26.23'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.53%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.814.98'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -0.46%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -4.10%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.71'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.42'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -3.52%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.38'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -4.93%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.88'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0306'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.62%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '328.76'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +5.61%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.33'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +9.89%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -3.45%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -4.92%  '
